# Auto-generated edit script applying the Famfrit_Profits diff
# Updates currentAveragePrice/currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ columns (H-N)
# across the 8 crafting-job worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 2469.3333
$ws.Range("I88").Value = 1003
$ws.Range("J88").Value = 2762.6
$ws.Range("K88").Value = 1003
$ws.Range("L88").Value = 2762.6
$ws.Range("M88").Value = -597
$ws.Range("N88").Value = -3574.6
$ws.Range("H91").Value = 2469.3333
$ws.Range("I91").Value = 1003
$ws.Range("J91").Value = 2762.6
$ws.Range("K91").Value = 1003
$ws.Range("L91").Value = 2762.6
$ws.Range("M91").Value = 401
$ws.Range("N91").Value = -5570.6
$ws.Range("H106").Value = 1993.9
$ws.Range("I106").Value = 1882.1111
$ws.Range("J106").Value = 3000
$ws.Range("K106").Value = 1882.1111
$ws.Range("L106").Value = 3000
$ws.Range("M106").Value = -1251.1111
$ws.Range("N106").Value = -4262
$ws.Range("H116").Value = 11021.8
$ws.Range("I116").Value = 11000.571
$ws.Range("J116").Value = 11071.333
$ws.Range("K116").Value = 11000.571
$ws.Range("L116").Value = 11071.333
$ws.Range("M116").Value = -7558.571
$ws.Range("N116").Value = -17955.333
$ws.Range("H132").Value = 7276.5557
$ws.Range("I132").Value = 7557.5293
$ws.Range("K132").Value = 22672.5879
$ws.Range("M132").Value = -20142.5879
$ws.Range("H135").Value = 20834488
$ws.Range("I135").Value = 1236
$ws.Range("K135").Value = 11124
$ws.Range("M135").Value = -8589

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8395.482
$ws.Range("I32").Value = 8284.477000000001
$ws.Range("J32").Value = 8686.875
$ws.Range("K32").Value = 8284.477000000001
$ws.Range("L32").Value = 8686.875
$ws.Range("M32").Value = -7997.477000000001
$ws.Range("N32").Value = -9260.875
$ws.Range("H63").Value = 2361
$ws.Range("I63").Value = 2019.8
$ws.Range("K63").Value = 2019.8
$ws.Range("M63").Value = -1333.8
$ws.Range("H66").Value = 2361
$ws.Range("I66").Value = 2019.8
$ws.Range("K66").Value = 10099
$ws.Range("M66").Value = -6667
$ws.Range("H76").Value = 250188
$ws.Range("J76").Value = 250188
$ws.Range("L76").Value = 250188
$ws.Range("N76").Value = -250864
$ws.Range("H79").Value = 250188
$ws.Range("J79").Value = 250188
$ws.Range("L79").Value = 250188
$ws.Range("N79").Value = -252528
$ws.Range("H88").Value = 12218.556
$ws.Range("I88").Value = 25663
$ws.Range("J88").Value = 1463
$ws.Range("K88").Value = 25663
$ws.Range("L88").Value = 1463
$ws.Range("M88").Value = -25257
$ws.Range("N88").Value = -2275
$ws.Range("H91").Value = 12218.556
$ws.Range("I91").Value = 25663
$ws.Range("J91").Value = 1463
$ws.Range("K91").Value = 25663
$ws.Range("L91").Value = 1463
$ws.Range("M91").Value = -24259
$ws.Range("N91").Value = -4271
$ws.Range("H97").Value = 2497.375
$ws.Range("I97").Value = 2629
$ws.Range("J97").Value = 1049.5
$ws.Range("K97").Value = 2629
$ws.Range("L97").Value = 1049.5
$ws.Range("M97").Value = -2133
$ws.Range("N97").Value = -2041.5
$ws.Range("H102").Value = 287878.16
$ws.Range("I102").Value = 287878.16
$ws.Range("K102").Value = 287878.16
$ws.Range("M102").Value = -286256.16
$ws.Range("H103").Value = 70000
$ws.Range("I103").Value = 70000
$ws.Range("K103").Value = 70000
$ws.Range("M103").Value = -68828
$ws.Range("H132").Value = 30355868
$ws.Range("I132").Value = 12553.038
$ws.Range("J132").Value = 143059620
$ws.Range("K132").Value = 37659.114
$ws.Range("L132").Value = 429178860
$ws.Range("M132").Value = -35129.114
$ws.Range("N132").Value = -429183920

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2397
$ws.Range("I20").Value = 1832.1578
$ws.Range("J20").Value = 3291.3333
$ws.Range("K20").Value = 1832.1578
$ws.Range("L20").Value = 3291.3333
$ws.Range("M20").Value = -1585.1578
$ws.Range("N20").Value = -3785.3333
$ws.Range("H82").Value = 23363.572
$ws.Range("I82").Value = 5146
$ws.Range("J82").Value = 37026.75
$ws.Range("K82").Value = 5146
$ws.Range("L82").Value = 37026.75
$ws.Range("M82").Value = -4763
$ws.Range("N82").Value = -37792.75
$ws.Range("H85").Value = 23363.572
$ws.Range("I85").Value = 5146
$ws.Range("J85").Value = 37026.75
$ws.Range("K85").Value = 5146
$ws.Range("L85").Value = 37026.75
$ws.Range("M85").Value = -3820
$ws.Range("N85").Value = -39678.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 13892011
$ws.Range("I31").Value = 2574.9375
$ws.Range("K31").Value = 2574.9375
$ws.Range("M31").Value = -2279.9375
$ws.Range("H34").Value = 13892011
$ws.Range("I34").Value = 2574.9375
$ws.Range("K34").Value = 2574.9375
$ws.Range("M34").Value = -2372.9375
$ws.Range("H99").Value = 1658
$ws.Range("I99").Value = 1658
$ws.Range("K99").Value = 1658
$ws.Range("M99").Value = -160
$ws.Range("H105").Value = 7115.5
$ws.Range("I105").Value = 938.3570999999999
$ws.Range("J105").Value = 50355.5
$ws.Range("K105").Value = 938.3570999999999
$ws.Range("L105").Value = 50355.5
$ws.Range("M105").Value = 808.6429000000001
$ws.Range("N105").Value = -53849.5
$ws.Range("H126").Value = 1658
$ws.Range("I126").Value = 1658
$ws.Range("K126").Value = 4974
$ws.Range("M126").Value = -2504
$ws.Range("H134").Value = 2684.1428
$ws.Range("I134").Value = 1773.0625
$ws.Range("J134").Value = 5599.6
$ws.Range("K134").Value = 5319.1875
$ws.Range("L134").Value = 16798.8
$ws.Range("M134").Value = -2784.1875
$ws.Range("N134").Value = -21868.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 1627
$ws.Range("I92").Value = 1205.8
$ws.Range("J92").Value = 1978
$ws.Range("K92").Value = 3617.4
$ws.Range("L92").Value = 5934
$ws.Range("M92").Value = -2369.4
$ws.Range("N92").Value = -8430
$ws.Range("H139").Value = 2561.6667
$ws.Range("I139").Value = 2576.2666
$ws.Range("J139").Value = 2488.6667
$ws.Range("K139").Value = 7728.7998
$ws.Range("L139").Value = 7466.000100000001
$ws.Range("M139").Value = -2588.7998
$ws.Range("N139").Value = -17746.0001
$ws.Range("H140").Value = 1587.9048
$ws.Range("I140").Value = 1281.1765
$ws.Range("K140").Value = 3843.5295
$ws.Range("M140").Value = 1336.4705

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 997.62964
$ws.Range("I97").Value = 990.95
$ws.Range("K97").Value = 990.95
$ws.Range("M97").Value = -494.95
$ws.Range("H102").Value = 8842.571
$ws.Range("I102").Value = 2479.6
$ws.Range("K102").Value = 2479.6
$ws.Range("M102").Value = -857.5999999999999
$ws.Range("H122").Value = 41668348
$ws.Range("I122").Value = 1586.5
$ws.Range("K122").Value = 4759.5
$ws.Range("M122").Value = -2309.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3214.4614
$ws.Range("I22").Value = 2736.6667
$ws.Range("J22").Value = 3624
$ws.Range("K22").Value = 2736.6667
$ws.Range("L22").Value = 3624
$ws.Range("M22").Value = -2441.6667
$ws.Range("N22").Value = -4214
$ws.Range("H27").Value = 3214.4614
$ws.Range("I27").Value = 2736.6667
$ws.Range("J27").Value = 3624
$ws.Range("K27").Value = 2736.6667
$ws.Range("L27").Value = 3624
$ws.Range("M27").Value = -2629.6667
$ws.Range("N27").Value = -3838
$ws.Range("H31").Value = 4993
$ws.Range("I31").Value = 15
$ws.Range("J31").Value = 7979.8
$ws.Range("K31").Value = 15
$ws.Range("L31").Value = 7979.8
$ws.Range("M31").Value = 233
$ws.Range("N31").Value = -8475.799999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = ""
$ws.Range("N12").Value = ""
$ws.Range("H51").Value = 28747
$ws.Range("I51").Value = 27500
$ws.Range("K51").Value = 27500
$ws.Range("M51").Value = -26990
$ws.Range("H62").Value = 7600.143
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 7600.143
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 7600.143
$ws.Range("M62").Value = ""
$ws.Range("N62").Value = -8848.143
$ws.Range("H65").Value = 7600.143
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 7600.143
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 38000.715
$ws.Range("M65").Value = ""
$ws.Range("N65").Value = -44240.715
$ws.Range("H132").Value = 2505.7942
$ws.Range("I132").Value = 1962.2084
$ws.Range("J132").Value = 3810.4
$ws.Range("K132").Value = 5886.6252
$ws.Range("L132").Value = 11431.2
$ws.Range("M132").Value = -3356.6252
$ws.Range("N132").Value = -16491.2
$ws.Range("H136").Value = 11233.08
$ws.Range("I136").Value = 4987.6
$ws.Range("K136").Value = 14962.8
$ws.Range("M136").Value = -12412.8
